$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.63409411907196
$ws.Range("B1").Value = 1.7994225025177
$ws.Range("C1").Value = 2.132314682006836
$ws.Range("D1").Value = 3.477210998535156
$ws.Range("E1").Value = 3.21173620223999
